# "changed the au sheet" / "put the inflation adjusted gold data back in
# column c because it was giving me issues"
#
# Au sheet: move the existing column C ("Primary commodity price old" /
# oil-adjusted mine-production numbers) out to a new column Q, and pull the
# inflation-adjusted price series that was living in column M back into
# column C (with a fresh header string), including the M21 style-of-growth
# formula for the final year.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Au")

# ------------------------------------------------------------------
# 1. Preserve old column C (header text + data) into new column Q first,
#    before we overwrite C.
# ------------------------------------------------------------------
$ws.Range("Q1").Value2 = $ws.Range("C1").Value2
$ws.Range("Q2").Value2 = $ws.Range("C2").Value2

for ($r = 3; $r -le 21; $r++) {
    $oldC = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 17).Value2 = $oldC
}

# ------------------------------------------------------------------
# 2. New column C header ("Primary commodity price " - trailing space,
#    distinct string from the old "Primary commodity price old").
# ------------------------------------------------------------------
$ws.Range("C1").Value2 = "Primary commodity price "

# ------------------------------------------------------------------
# 3. Pull the inflation-adjusted price series (column M) back into C.
# ------------------------------------------------------------------
for ($r = 3; $r -le 20; $r++) {
    $mVal = $ws.Cells.Item($r, 13).Value2
    $ws.Cells.Item($r, 3).Value2 = $mVal
}

# Row 21 in M uses a growth formula off the prior year; mirror it in C.
$ws.Range("C21").Formula = "=C20+(C20*0.18)"

# ------------------------------------------------------------------
# 4. Selection moves to C1 after the rework.
# ------------------------------------------------------------------
$ws.Range("C1").Select()
